$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4799959659576416
$ws.Range("B1").Value = 1.662782073020935
$ws.Range("C1").Value = 3.401604413986206
$ws.Range("D1").Value = 3.677681922912598
$ws.Range("E1").Value = 1.471039533615112
